$wb = $excel.ActiveWorkbook

# Sheet "Hoja1" contains the long conversion text in A1
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.41 = 29818.52 pesos`n✅ 29818.52 pesos = 7.36 = 968.37 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# Sheet "tasas" contains the N10/O10/N12/O12 numeric updates
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 135
$ws2.Range("O10").Value = 4025.5

$ws2.Range("N12").Value = 4049.5
$ws2.Range("O12").Value = 131.51
